# Update CO2 price results per sheet (year) with new server results.
$wb = $excel.ActiveWorkbook

$sheetNames = @("2025", "2030", "2040", "2045", "2050")
$newValues  = @(71.25,   243.75,  443.75,  443.75,  443.75)

for ($i = 0; $i -lt $sheetNames.Length; $i++) {
    $ws = $wb.Worksheets.Item($sheetNames[$i])
    $ws.Range("A2").Value = $newValues[$i]
}
